$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new contact row (row 4) that was previously empty
$ws.Range("A4").Value = "Sebastian Romero"
$ws.Range("B4").Value = 573138793438
$ws.Range("C4").Value = "Joven"

# Update the active cell selection to C5 (as in the source workbook)
$ws.Range("C5").Select()
